$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 41 / 42 swap: Maker <-> Kaspa (including B "Coin" and C "Link" columns) ---
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.136"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.390.30"
$ws.Range("E42").Value = "  +5.99%  "

# --- Price (column D) updates ---
$ws.Range("D2").Value = "70.733.06"
$ws.Range("D3").Value = "3.583.37"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.23"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.45"
$ws.Range("D7").Value = "3.572.44"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.216"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.652"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.21"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000323"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.57"
$ws.Range("D15").Value = "4.154.92"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.58"
$ws.Range("D17").Value = "70.719.27"
$ws.Range("D18").Value = "3.605.10"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "573.76"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.40"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.71"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.93"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.17"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.38"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.95"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.15"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.31"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.35"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.90"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.114"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "562.74"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.418"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.67"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.38"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.58"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0448"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.97"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.33"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.45"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  +17.43%  "
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("E13").Value = "  +5.81%  "
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("E19").Value = "  +16.28%  "
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  -2.19%  "
$ws.Range("E23").Value = "  -9.19%  "
$ws.Range("E24").Value = "  +5.58%  "
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("E29").Value = "  -2.56%  "
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("E31").Value = "  -5.13%  "
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("E35").Value = "  +4.60%  "
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("E38").Value = "  -3.76%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("E44").Value = "  -2.50%  "
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("E47").Value = "  -3.81%  "
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("E51").Value = "  -7.36%  "
